$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 0.001
$ws.Range("K5").Value = 713
$ws.Range("L5").Value = 0.001426
